# ------------------------------------------------------------------
# Fix-banred-backsub, fix axiales, FIX EA, FINAL PROGRAM
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# ===================================================================
# Sheet "Elementos" - element connectivity table (DIM, Conectiv_i, Conectiv_j)
# ===================================================================
$wsElem = $wb.Worksheets.Item("Elementos")

# Row 4 : i/j changed
$wsElem.Range("B4").Value = 2
$wsElem.Range("C4").Value = 3

# Row 6 : i/j changed
$wsElem.Range("B6").Value = 0
$wsElem.Range("C6").Value = 4

# Row 7 : i/j changed
$wsElem.Range("B7").Value = 1
$wsElem.Range("C7").Value = 5

# Rows 8-21 : new members (DIM already carries style, i/j are new values)
$wsElem.Range("A8").Value = 1
$wsElem.Range("B8").Value = 2
$wsElem.Range("C8").Value = 6

$wsElem.Range("A9").Value = 1
$wsElem.Range("B9").Value = 3
$wsElem.Range("C9").Value = 7

$wsElem.Range("A10").Value = 1
$wsElem.Range("B10").Value = 4
$wsElem.Range("C10").Value = 5

$wsElem.Range("A11").Value = 1
$wsElem.Range("B11").Value = 5
$wsElem.Range("C11").Value = 6

$wsElem.Range("A12").Value = 1
$wsElem.Range("B12").Value = 6
$wsElem.Range("C12").Value = 7

$wsElem.Range("A13").Value = 1
$wsElem.Range("B13").Value = 4
$wsElem.Range("C13").Value = 7

$wsElem.Range("A14").Value = 1
$wsElem.Range("B14").Value = 0
$wsElem.Range("C14").Value = 8

$wsElem.Range("A15").Value = 1
$wsElem.Range("B15").Value = 1
$wsElem.Range("C15").Value = 8

$wsElem.Range("A16").Value = 1
$wsElem.Range("B16").Value = 2
$wsElem.Range("C16").Value = 8

$wsElem.Range("A17").Value = 1
$wsElem.Range("B17").Value = 3
$wsElem.Range("C17").Value = 8

$wsElem.Range("A18").Value = 1
$wsElem.Range("B18").Value = 4
$wsElem.Range("C18").Value = 8

$wsElem.Range("A19").Value = 1
$wsElem.Range("B19").Value = 5
$wsElem.Range("C19").Value = 8

$wsElem.Range("A20").Value = 1
$wsElem.Range("B20").Value = 6
$wsElem.Range("C20").Value = 8

$wsElem.Range("A21").Value = 1
$wsElem.Range("B21").Value = 7
$wsElem.Range("C21").Value = 8

$wsElem.Activate()
$wsElem.Range("G17").Select()

# ===================================================================
# Sheet "Nodos" - node coordinates / loads / restraints
# ===================================================================
$wsNodos = $wb.Worksheets.Item("Nodos")

# Copy an already-formatted (bordered) empty cell so new cells pick up
# the same style (s="1") without inventing new style entries.
$fmtSrc = $wsNodos.Range("A6")
$fmtSrc.Copy()
$wsNodos.Range("D2:I5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2 : D2..I2 become 0 (values already 0, formatting now applied above)
$wsNodos.Range("D2").Value = 0
$wsNodos.Range("E2").Value = 0
$wsNodos.Range("F2").Value = 0
$wsNodos.Range("G2").Value = 0
$wsNodos.Range("H2").Value = 0
$wsNodos.Range("I2").Value = 0

# Row 3
$wsNodos.Range("B3").Value = 1
$wsNodos.Range("D3").Value = 0
$wsNodos.Range("E3").Value = 0
$wsNodos.Range("F3").Value = 0
$wsNodos.Range("G3").Value = 0
$wsNodos.Range("H3").Value = 0
$wsNodos.Range("I3").Value = 0

# Row 4
$wsNodos.Range("A4").Value = 1
$wsNodos.Range("B4").Value = 1
$wsNodos.Range("D4").Value = 0
$wsNodos.Range("E4").Value = 0
$wsNodos.Range("F4").Value = 0
$wsNodos.Range("G4").Value = 0
$wsNodos.Range("H4").Value = 0
$wsNodos.Range("I4").Value = 0

# Row 5
$wsNodos.Range("A5").Value = 1
$wsNodos.Range("B5").Value = 0
$wsNodos.Range("C5").Value = 0
$wsNodos.Range("D5").Value = 0
$wsNodos.Range("E5").Value = 0
$wsNodos.Range("F5").Value = 0
$wsNodos.Range("G5").Value = 0
$wsNodos.Range("H5").Value = 0
$wsNodos.Range("I5").Value = 0

# Row 6 (new load case data - load of -10 kN in FZ, fully restrained)
$wsNodos.Range("A6").Value = 0
$wsNodos.Range("B6").Value = 0
$wsNodos.Range("C6").Value = 1
$wsNodos.Range("D6").Value = 0
$wsNodos.Range("E6").Value = 0
$wsNodos.Range("F6").Value = -10
$wsNodos.Range("G6").Value = 1
$wsNodos.Range("H6").Value = 1
$wsNodos.Range("I6").Value = 1

# Row 7
$wsNodos.Range("A7").Value = 0
$wsNodos.Range("B7").Value = 1
$wsNodos.Range("C7").Value = 1
$wsNodos.Range("D7").Value = 0
$wsNodos.Range("E7").Value = 0
$wsNodos.Range("F7").Value = -10
$wsNodos.Range("G7").Value = 1
$wsNodos.Range("H7").Value = 1
$wsNodos.Range("I7").Value = 1

# Row 8
$wsNodos.Range("A8").Value = 1
$wsNodos.Range("B8").Value = 1
$wsNodos.Range("C8").Value = 1
$wsNodos.Range("D8").Value = 0
$wsNodos.Range("E8").Value = 0
$wsNodos.Range("F8").Value = -10
$wsNodos.Range("G8").Value = 1
$wsNodos.Range("H8").Value = 1
$wsNodos.Range("I8").Value = 1

# Row 9
$wsNodos.Range("A9").Value = 1
$wsNodos.Range("B9").Value = 0
$wsNodos.Range("C9").Value = 1
$wsNodos.Range("D9").Value = 0
$wsNodos.Range("E9").Value = 0
$wsNodos.Range("F9").Value = -10
$wsNodos.Range("G9").Value = 1
$wsNodos.Range("H9").Value = 1
$wsNodos.Range("I9").Value = 1

# Row 10
$wsNodos.Range("A10").Value = 0.5
$wsNodos.Range("B10").Value = 0.5
$wsNodos.Range("C10").Value = 1
$wsNodos.Range("D10").Value = 0
$wsNodos.Range("E10").Value = 0
$wsNodos.Range("F10").Value = -10
$wsNodos.Range("G10").Value = 1
$wsNodos.Range("H10").Value = 1
$wsNodos.Range("I10").Value = 1

$wsNodos.Range("C10").Select()

# ===================================================================
# Sheet "Datos" - no data change (only shared-string reindex, automatic)
# ===================================================================
$wsDatos = $wb.Worksheets.Item("Datos")

# ===================================================================
# Sheet "Props" - cross-section properties (EA column reworked into
# separate Area / Young's-modulus columns, plus helper labels)
# ===================================================================
$wsProps = $wb.Worksheets.Item("Props")

$wsProps.Range("C1").Value = "A (cm^2)"
$wsProps.Range("E1").Value = "A en cm^2"
$wsProps.Range("F1").Value = "AE en 10⁵N (100 kN)."

$wsProps.Range("B2").Value = 100
$wsProps.Range("C2").Value = 10
$wsProps.Range("E2").Value = "E en GPA (kN/mm^2)"

$wsProps.Range("B3").Value = 100
$wsProps.Range("C3").Value = 10

$wsProps.Range("B4").Value = 100
$wsProps.Range("C4").Value = 10

$wsProps.Activate()
$wsProps.Range("F9").Select()
